# Apply the "Adding 2D plotting script with gif!" edit:
#  - Add a new summary table (rows 9-15) to the "Agility" sheet, mirroring the
#    existing table but for the "10% flooded" configuration.
#  - Make "Agility" the active sheet/tab (it was "Flood %" before).
#  - Update the selection/active-cell on each sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Agility")
$ws2 = $wb.Worksheets.Item("Flood %")

# ---- New table on the Agility sheet (rows 9-13), row 14 blank, row 15 label ----

# Row 9: header values (same agility thresholds as row 1), scientific formatted
$ws1.Range("A9").Value = "Config"
$ws1.Range("B9").NumberFormat = "0.00E+00"
$ws1.Range("C9").NumberFormat = "0.00E+00"
$ws1.Range("D9").NumberFormat = "0.00E+00"
$ws1.Range("B9").Value = 0.0000004
$ws1.Range("C9").Value = 0.00004
$ws1.Range("D9").Value = 0.004
$ws1.Range("E9").Value = "# of unique alt outliers"

# Row 10: Nadir
$ws1.Range("A10").Value = "Nadir"
$ws1.Range("B10").Value = 20
$ws1.Range("C10").Value = 20
$ws1.Range("D10").Value = 20
$ws1.Range("E10").Value = 98

# Row 11: Agile
$ws1.Range("A11").Value = "Agile"
$ws1.Range("B11").Value = 16
$ws1.Range("C11").Value = 36
$ws1.Range("D11").Value = 41
$ws1.Range("E11").Value = 98

# Row 12: No preplan 3D
$ws1.Range("A12").Value = "No preplan 3D"
$ws1.Range("B12").Value = 48
$ws1.Range("C12").Value = 58
$ws1.Range("D12").Value = 58
$ws1.Range("E12").Value = 98

# Row 13: Preplan 3D
$ws1.Range("A13").Value = "Preplan 3D"
$ws1.Range("B13").Value = 38
$ws1.Range("C13").Value = 37
$ws1.Range("D13").Value = 39
$ws1.Range("E13").Value = 98

# Row 14 intentionally left blank (spacer row)

# Row 15: label for the new table
$ws1.Range("A15").Value = "All with 10% flooded"

# ---- Selection / active sheet updates ----

# Flood % sheet: move selection to D6, no longer the active tab
[void]$ws2.Range("D6").Select()

# Agility sheet: move selection to B13, becomes the active tab (selected last)
[void]$ws1.Range("B13").Select()
